$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'41.561.08"
$ws.Range("E2").Value = "'  -1.18%  "
$ws.Range("D3").Value = "'2.212.68"
$ws.Range("E3").Value = "'  -0.69%  "
$ws.Range("E4").Value = "'  +0.08%  "
$ws.Range("D5").Value = "'253.56"
$ws.Range("E5").Value = "'  +4.49%  "
$ws.Range("D6").Value = "'0.629"
$ws.Range("E6").Value = "'  +0.48%  "
$ws.Range("D7").Value = "'69.91"
$ws.Range("E7").Value = "'  +2.82%  "
$ws.Range("E8").Value = "'  +0.07%  "
$ws.Range("D9").Value = "'0.596"
$ws.Range("E9").Value = "'  +9.56%  "
$ws.Range("D10").Value = "'39.54"
$ws.Range("E10").Value = "'  +11.74%  "
$ws.Range("D11").Value = "'0.0961"
$ws.Range("E11").Value = "'  -0.09%  "
$ws.Range("D12").Value = "'58.50"
$ws.Range("E12").Value = "'  +0.29%  "
$ws.Range("D13").Value = "'7.23"
$ws.Range("E13").Value = "'  +8.18%  "
$ws.Range("E14").Value = "'  -0.01%  "
$ws.Range("D15").Value = "'2.545.18"
$ws.Range("E15").Value = "'  -0.64%  "
$ws.Range("D16").Value = "'0.894"
$ws.Range("E16").Value = "'  +5.17%  "
$ws.Range("D17").Value = "'14.92"
$ws.Range("E17").Value = "'  +0.98%  "
$ws.Range("D18").Value = "'2.230.85"
$ws.Range("E18").Value = "'  +0.07%  "
$ws.Range("D19").Value = "'41.581.37"
$ws.Range("E19").Value = "'  -0.90%  "
$ws.Range("D20").Value = "'0.0₃0963"
$ws.Range("E20").Value = "'  +0.83%  "
$ws.Range("E21").Value = "'  +2.73%  "
$ws.Range("D22").Value = "'72.44"
$ws.Range("E22").Value = "'  -0.18%  "
$ws.Range("D23").Value = "'233.95"
$ws.Range("E23").Value = "'  +0.02%  "
$ws.Range("E24").Value = "'  +1.04%  "
$ws.Range("D25").Value = "'3.99"
$ws.Range("E25").Value = "'  +11.32%  "
$ws.Range("D26").Value = "'11.90"
$ws.Range("E26").Value = "'  +20.66%  "
$ws.Range("E28").Value = "'  +4.30%  "
$ws.Range("E29").Value = "'  -1.74%  "
$ws.Range("D30").Value = "'172.34"
$ws.Range("E30").Value = "'  +0.46%  "
$ws.Range("D31").Value = "'20.69"
$ws.Range("E31").Value = "'  +1.62%  "
$ws.Range("E32").Value = "'  +2.59%  "
$ws.Range("E33").Value = "'  +9.24%  "
$ws.Range("D34").Value = "'0.124"
$ws.Range("E34").Value = "'  -0.90%  "
$ws.Range("D35").Value = "'0.0742"
$ws.Range("E35").Value = "'  +4.37%  "
$ws.Range("B36").Value = "'InjectiveProtocol"
$ws.Range("C36").Value = "'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D36").Value = "'26.45"
$ws.Range("E36").Value = "'  +15.30%  "
$ws.Range("B37").Value = "'Filecoin"
$ws.Range("C37").Value = "'https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D37").Value = "'4.70"
$ws.Range("E37").Value = "'  +0.97%  "
$ws.Range("D38").Value = "'4.01"
$ws.Range("E38").Value = "'  +3.42%  "
$ws.Range("E39").Value = "'  +7.58%  "
$ws.Range("E40").Value = "'  -0.95%  "
$ws.Range("D41").Value = "'5.89"
$ws.Range("E41").Value = "'  +1.50%  "
$ws.Range("D42").Value = "'12.24"
$ws.Range("E42").Value = "'  +25.44%  "
$ws.Range("D43").Value = "'64.72"
$ws.Range("E43").Value = "'  -3.13%  "
$ws.Range("B44").Value = "'Algorand"
$ws.Range("C44").Value = "'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D44").Value = "'0.203"
$ws.Range("E44").Value = "'  +6.24%  "
$ws.Range("B45").Value = "'FTXToken"
$ws.Range("C45").Value = "'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D45").Value = "'4.90"
$ws.Range("E45").Value = "'  -0.32%  "
$ws.Range("D46").Value = "'8.81"
$ws.Range("E46").Value = "'  -2.21%  "
$ws.Range("E47").Value = "'  +1.24%  "
$ws.Range("E48").Value = "'  +0.38%  "
$ws.Range("E49").Value = "'  +6.02%  "
$ws.Range("D50").Value = "'4.51"
$ws.Range("E50").Value = "'  -2.46%  "
$ws.Range("D51").Value = "'2.40"
$ws.Range("E51").Value = "'  +4.68%  "
